$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06373566666666666
$ws.Range("H2").Value = 0.191207
$ws.Range("I2").Value = 0.01058875298517695
$ws.Range("J2").Value = 0.01058875298517695
$ws.Range("M2").Value = 34.53319033333333
$ws.Range("N2").Value = 103.599571
$ws.Range("O2").Value = 0.2461870921144496
$ws.Range("P2").Value = 0.2461870921144496
$ws.Range("Q2").Value = 2.200995908021889
$ws.Range("R2").Value = 19.808963172197
$ws.Range("S2").Value = 0.002606814306538912
$ws.Range("T2").Value = 0.002606814306538912
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06373566666666666
$ws.Range("H3").Value = 0.191207
$ws.Range("I3").Value = 0.01058875298517695
$ws.Range("J3").Value = 0.01058875298517695
$ws.Range("O3").Value = 0.575843103803214
$ws.Range("P3").Value = 0.575843103803214
$ws.Range("Q3").Value = 5.148232201159777
$ws.Range("R3").Value = 46.33408981043799
$ws.Range("S3").Value = 0.006097460384389845
$ws.Range("T3").Value = 0.006097460384389845
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06373566666666666
$ws.Range("H4").Value = 0.191207
$ws.Range("I4").Value = 0.01058875298517695
$ws.Range("J4").Value = 0.01058875298517695
$ws.Range("O4").Value = 0.1779698040823365
$ws.Range("P4").Value = 0.1779698040823364
$ws.Range("Q4").Value = 1.591110269723556
$ws.Range("R4").Value = 14.319992427512
$ws.Range("S4").Value = 0.001884478294248198
$ws.Range("T4").Value = 0.001884478294248197
$ws.Range("I5").Value = 0.2961697031425515
$ws.Range("J5").Value = 0.2961697031425515
$ws.Range("M5").Value = 34.53319033333333
$ws.Range("N5").Value = 103.599571
$ws.Range("O5").Value = 0.2461870921144496
$ws.Range("P5").Value = 0.2461870921144496
$ws.Range("Q5").Value = 61.56232991829677
$ws.Range("R5").Value = 554.0609692646709
$ws.Range("S5").Value = 0.07291315798906453
$ws.Range("T5").Value = 0.07291315798906452
$ws.Range("I6").Value = 0.2961697031425515
$ws.Range("J6").Value = 0.2961697031425515
$ws.Range("O6").Value = 0.575843103803214
$ws.Range("P6").Value = 0.575843103803214
$ws.Range("S6").Value = 0.1705472811100834
$ws.Range("T6").Value = 0.1705472811100834
$ws.Range("I7").Value = 0.2961697031425515
$ws.Range("J7").Value = 0.2961697031425515
$ws.Range("O7").Value = 0.1779698040823365
$ws.Range("P7").Value = 0.1779698040823364
$ws.Range("S7").Value = 0.05270926404340364
$ws.Range("T7").Value = 0.05270926404340363
$ws.Range("I8").Value = 0.6932415438722715
$ws.Range("J8").Value = 0.6932415438722715
$ws.Range("M8").Value = 34.53319033333333
$ws.Range("N8").Value = 103.599571
$ws.Range("O8").Value = 0.2461870921144496
$ws.Range("P8").Value = 0.2461870921144496
$ws.Range("Q8").Value = 144.0983469412897
$ws.Range("R8").Value = 1296.885122471608
$ws.Range("S8").Value = 0.1706671198188462
$ws.Range("T8").Value = 0.1706671198188462
$ws.Range("I9").Value = 0.6932415438722715
$ws.Range("J9").Value = 0.6932415438722715
$ws.Range("O9").Value = 0.575843103803214
$ws.Range("P9").Value = 0.575843103803214
$ws.Range("S9").Value = 0.3991983623087408
$ws.Range("T9").Value = 0.3991983623087408
$ws.Range("I10").Value = 0.6932415438722715
$ws.Range("J10").Value = 0.6932415438722715
$ws.Range("O10").Value = 0.1779698040823365
$ws.Range("P10").Value = 0.1779698040823364
$ws.Range("S10").Value = 0.1233760617446846
$ws.Range("T10").Value = 0.1233760617446846
